# "Generate Report for Handback" - refresh the handback-status report's
# timestamp columns (Latest HO Xliff Generate Date / Correspond Handoff
# Datetime / Correspond Handback DateTime) to the latest run's values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-29 19:11:08"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime".
$wsZhCn.Range("H2").Value = "2016-08-29 19:10:57"
$wsZhCn.Range("K2").Value = "2016-08-29 19:11:30"

# de-de sheet: "Correspond Handback DateTime".
$wsDeDe.Range("K2").Value = "2016-08-29 19:11:37"
